# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview + per-locale "Status" cells move from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - zh-cn sheet: Latest Handback DateTime refreshed, stale "Error Detail"
#    messages cleared now that the handback is in sync
#  - de-de sheet: handback completed - Latest Target File / Latest Handback
#    File / Latest Handback DateTime populated, with a new hyperlink on the
#    Latest Target File cell (mirroring the zh-cn sheet)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$statusDone = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = $statusDone
$overview.Range("F2").Value = $statusDone
$overview.Range("E3").Value = $statusDone
$overview.Range("F3").Value = $statusDone

$zhcn.Range("C2").Value = $statusDone
$zhcn.Range("C3").Value = $statusDone

$dede.Range("C2").Value = $statusDone
$dede.Range("C3").Value = $statusDone

# --- zh-cn: refresh Latest Handback DateTime, clear stale Error Detail ---
$zhcn.Range("K2").Value = "2016-09-07 06:53:53"
$zhcn.Range("K3").Value = "2016-09-07 06:53:53"

$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# --- de-de: handback completed, populate target/handback file + datetime ---
$dedeHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("J2").Value = $dedeHandbackFile
$dede.Range("J3").Value = $dedeHandbackFile

$dede.Range("K2").Value = "2016-09-07 06:54:04"
$dede.Range("K3").Value = "2016-09-07 06:54:04"

# Rebuild de-de hyperlinks in row order (A2, I2, A3, I3) so relationship ids
# line up the same way they already do on the zh-cn sheet.
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e329d65c1d2dd70e658bbf803e07ab06af65dd96/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/14e4cf8d03dede0905ab6131718f7fc246f37418/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e329d65c1d2dd70e658bbf803e07ab06af65dd96/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/14e4cf8d03dede0905ab6131718f7fc246f37418/e2e/a.md", "", "", "a.md")

# --- Column width tweaks (wider Status / narrower Error Detail columns) ---
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667

$zhcn.Columns.Item(3).ColumnWidth  = 29.16666667
$zhcn.Columns.Item(16).ColumnWidth = 12.83333333

$dede.Columns.Item(3).ColumnWidth  = 29.16666667
$dede.Columns.Item(10).ColumnWidth = 39.16666667

Write-Host "Handback report generated"
